# Update generated output numbers (commit: "Update gh-pages to output generated at 456a3b4")
# Column F ("想去人数" / interest count) changes on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1038
$wsExpo.Range("F3").Value = 59
$wsExpo.Range("F4").Value = 2313
$wsExpo.Range("F6").Value = 508

# Sheet "全部类型" (all types, combined listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1038
$wsAll.Range("F5").Value = 59
$wsAll.Range("F6").Value = 2313
$wsAll.Range("F8").Value = 508
